$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3 V 0.3")

# Fix E3: was stored as inline string "532900", should be a numeric value.
$ws.Range("E3").Value = 532900

# Add new row 4 with the breakout data.
$ws.Range("A4").Value = "12/06/2024 07:44:47"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "PAISALO"
$ws.Range("D4").Value = "Paisalo Digital Ltd"

# E4's bsecode is stored as text (not a number) in the target workbook, so
# force the cell to text format before assigning the numeric-looking string,
# otherwise Excel's automatic type inference would store it as a number.
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "532900"
$ws.Range("E4").Style = "Normal"

$ws.Range("F4").Value = 9.17
$ws.Range("G4").Value = 69.39
$ws.Range("H4").Value = 3678180
